$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("zh-cn")
$ws1.Range("E4").Value = "2016-03-18 12:31:38"
$ws1.Range("H4").Value = "2016-03-18 12:31:57"

$ws2 = $wb.Worksheets.Item("de-de")
$ws2.Range("E4").Value = "2016-03-18 12:31:41"
$ws2.Range("H4").Value = "2016-03-18 12:32:03"
